$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are swapped between row 2 and row 3
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "{0}2" -f $col
    $addr3 = "{0}3" -f $col

    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2

    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
